# Fruta / hortaliza, semanal
# Insert a new data row (row 510) into the "Pepino ensalada" sheet, pushing
# the former row 510 down to row 511. The new row duplicates the values
# that were in row 509 (same date/quality/price entry recorded again).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 510 (and below) down by one row.
$ws.Rows(510).Insert()

# Populate the newly inserted row 510 with its data.
$ws.Cells.Item(510, 1).Value = 10
$ws.Cells.Item(510, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(510, 3).Value = "La Araucanía"
$ws.Cells.Item(510, 4).Value = 44442
$ws.Cells.Item(510, 5).Value = 9
$ws.Cells.Item(510, 6).Value = 100112043
$ws.Cells.Item(510, 7).Value = "Pepino ensalada"
$ws.Cells.Item(510, 8).Value = "Sin especificar"
$ws.Cells.Item(510, 9).Value = "Primera"
$ws.Cells.Item(510, 10).Value = 100
$ws.Cells.Item(510, 11).Value = 18000
$ws.Cells.Item(510, 12).Value = 18000
$ws.Cells.Item(510, 13).Value = 18000
$ws.Cells.Item(510, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(510, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(510, 16).Value = 300
$ws.Cells.Item(510, 17).Value = 60
$ws.Cells.Item(510, 18).Value = "Hortaliza"
